$d = $word.ActiveDocument

# The document opens with a centered name paragraph ("Dheeraj Chand"), and the
# contact-info line that should follow it is missing (short-resume bug). Re-write
# that first paragraph's range as itself plus a new, centered contact-info
# paragraph immediately after it, using raw OOXML so the new paragraph/run gets
# no inherited (bold/large) character formatting from the name run.

$firstPara = $d.Paragraphs.First
$nameRange = $firstPara.Range

$contactText = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:jc w:val="center"/></w:pPr>
            <w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>Dheeraj Chand</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:jc w:val="center"/></w:pPr>
            <w:r><w:t>$contactText</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

[void]$nameRange.InsertXML($xml)
